$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 133.3333333333333
$ws.Range("B4").Value = 8362.283772317136
$ws.Range("B5").Value = 121802.6639999999
$ws.Range("B7").Value = 1421.18598771462
$ws.Range("B8").Value = 20217.59999999971
$ws.Range("B9").Value = 2865.882758648496
$ws.Range("B10").Value = 155763.1820955332
$ws.Range("B11").Value = 0.07527144161574192
$ws.Range("B12").Value = 0.263148338158344
$ws.Range("B13").Value = 0.3499999999999986
$ws.Range("B14").Value = 0.9987371867967868
$ws.Range("B15").Value = 0.9244807679322856
